$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Ridge)
$ws.Range("B2").Value = 0.9999960958032298
$ws.Range("C2").Value = 0.2895728968386183
$ws.Range("D2").Value = 0.7104231989646115
$ws.Range("E2").Value = 0.06481042641419571
$ws.Range("F2").Value = 27.64641295283707
$ws.Range("G2").Value = "{'solver': 'saga', 'alpha': 0.001}"
$ws.Range("H2").Value = 2.17

# Row 3 (Lasso)
$ws.Range("B3").Value = 0.9999923035685135
$ws.Range("C3").Value = 0.6435909691867521
$ws.Range("D3").Value = 0.3564013343817615
$ws.Range("E3").Value = 0.09099629998478549
$ws.Range("F3").Value = 19.5818346861296
$ws.Range("H3").Value = 1.12

# Row 4 (ElasticNet)
$ws.Range("B4").Value = 0.9999998292898571
$ws.Range("C4").Value = 0.2851446271829342
$ws.Range("D4").Value = 0.7148552021069229
$ws.Range("E4").Value = 0.01355216538312895
$ws.Range("F4").Value = 27.73244260013872
$ws.Range("G4").Value = "{'l1_ratio': 0.9, 'alpha': 0.001}"
$ws.Range("H4").Value = 2.85

# Row 5 (SVR)
$ws.Range("B5").Value = 0.9567770954856544
$ws.Range("D5").Value = 0.4016976419844881
$ws.Range("E5").Value = 6.819240442267581
$ws.Range("H5").Value = 2.27

# Row 6 (KNN Regressor)
$ws.Range("H6").Value = 1.58

# Row 7 (Decision Tree)
$ws.Range("C7").Value = -1.094332560472273
$ws.Range("D7").Value = 1.265617476236871
$ws.Range("F7").Value = 47.46810040372116
$ws.Range("H7").Value = 2.37

# Row 8 (PLSRegression)
$ws.Range("C8").Value = 0.2277983465155331
$ws.Range("D8").Value = 0.7703424069921301
$ws.Range("E8").Value = 1.414319454662002
$ws.Range("F8").Value = 28.82334589671152
$ws.Range("H8").Value = 0.54
